# Update dashboards - 2026-01-17
# Advance the "as of" dates by one day and roll the Q:U value windows
# forward for the affected rows (29, 30, 47, 48, 49, 50, 52) on the
# "Aguilar Prototype" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 ---
$ws.Range("N29").Value = 46038
$ws.Range("Q29").Value = 2.27
$ws.Range("R29").Value = 2.22
$ws.Range("S29").Value = 2.22
$ws.Range("T29").Value = 2.23
$ws.Range("U29").Value = 2.22

# --- Row 30 ---
$ws.Range("N30").Value = 46038
$ws.Range("Q30").Value = 2.33
$ws.Range("R30").Value = 2.29
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.3
$ws.Range("U30").Value = 2.29

# --- Row 47 (date only) ---
$ws.Range("N47").Value = 46037

# --- Row 48 ---
$ws.Range("N48").Value = 46037
$ws.Range("Q48").Value = 3.56
$ws.Range("R48").Value = 3.51
$ws.Range("S48").Value = 3.53
$ws.Range("T48").Value = 3.54
$ws.Range("U48").Value = 3.54

# --- Row 49 ---
$ws.Range("N49").Value = 46037
$ws.Range("Q49").Value = 3.77
$ws.Range("R49").Value = 3.72
$ws.Range("S49").Value = 3.75
$ws.Range("T49").Value = 3.77
$ws.Range("U49").Value = 3.75

# --- Row 50 ---
$ws.Range("N50").Value = 46037
$ws.Range("Q50").Value = 4.17
$ws.Range("R50").Value = 4.15
$ws.Range("S50").Value = 4.18
$ws.Range("T50").Value = 4.19
$ws.Range("U50").Value = 4.18

# --- Row 52 ---
$ws.Range("N52").Value = 46037
$ws.Range("Q52").Value = 5.82
$ws.Range("R52").Value = 5.83
$ws.Range("S52").Value = 5.87
$ws.Range("T52").Value = 5.89
$ws.Range("U52").Value = 5.88
